# Applies a permutation of the species-observation data rows (2-14) on the
# "Artfynd" sheet. Every row keeps its shared/common metadata columns
# (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) exactly
# where they are; only the per-species columns actually move between rows:
#   A  Id
#   B  Taxonsorteringsordning
#   D  Rodlistade
#   E  TaxonId
#   F  Artnamn
#   G  Vetenskapligt namn
#   H  Auktor
#   Q  Ost
#   R  Nord
#   AC Publik kommentar   (present on some rows only)
#   AF Bestamningsmetod   (present on some rows only, always empty)
#
# Mapping below: new row number -> old row number (i.e. the data that used
# to live in row $mapping[$r] now belongs in row $r).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 14

$mapping = @{
  2  = 3
  3  = 8
  4  = 12
  5  = 10
  6  = 11
  7  = 14
  8  = 2
  9  = 13
  10 = 5
  11 = 7
  12 = 6
  13 = 9
  14 = 4
}

# Columns that move as a simple, always-populated value block.
$simpleCols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# 1. Snapshot the simple columns for every source row up front (this is a
#    full permutation, so writing early could clobber a value that is still
#    needed as a source for another destination row later on).
$simpleSnapshots = @{}
foreach ($col in $simpleCols) {
  $colSnap = @{}
  for ($r = $firstRow; $r -le $lastRow; $r++) {
    $colSnap[$r] = $ws.Range($col + $r).Value()
  }
  $simpleSnapshots[$col] = $colSnap
}

# Snapshot AC (public comment) and AF (determination-method) too, recording
# whether each source row actually had that cell at all.
$acSnapshots = @{}
$afSnapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $acSnapshots[$r] = $ws.Range("AC" + $r).Value()
  $afSnapshots[$r] = $ws.Range("AF" + $r).Value()
}

# 2. Write the simple columns back out to their destination rows.
foreach ($col in $simpleCols) {
  $colSnap = $simpleSnapshots[$col]
  for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $ws.Range($col + $r).Value = $colSnap[$srcRow]
  }
}

# 3. Write AC / AF, preserving "cell absent" vs "cell present but empty" vs
#    "cell has text" for each destination row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $srcRow = $mapping[$r]

  $acVal = $acSnapshots[$srcRow]
  if ($acVal -eq $null) {
    $ws.Range("AC" + $r).ClearContents()
  } else {
    $ws.Range("AC" + $r).Value = $acVal
    $ws.Range("AC" + $r).Style = "Normal"
  }

  $afVal = $afSnapshots[$srcRow]
  if ($afVal -eq $null) {
    $ws.Range("AF" + $r).ClearContents()
  } else {
    if ($afVal -eq "") {
      $ws.Range("AF" + $r).Value = "'"
    } else {
      $ws.Range("AF" + $r).Value = $afVal
    }
    $ws.Range("AF" + $r).Style = "Normal"
  }
}
